# 2ИСИП-722_ТерВер.xlsx — update attendance/score grid and selection state
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header row 3: clear the point values that used to sit in H3:J3
$ws.Range("H3:J3").ClearContents()

# Every zero-valued score cell in C4:J32 becomes 2 (21 contiguous row-blocks
# covering the 85 individual cells touched by the diff). The dependent
# SUM() formulas in column L recalc automatically.
$scoreRanges = $ws.Range("G4:J4,H6:J6,C7:J7,G8:J8,J9,I11:J11,I14:J14,I15:J15,C17:J17,G18:J18,G19:J19,C21:J21,E23:J23,J25,C26:J26,C28:J28,J29,C30:J30,J31,E32,I32")
foreach ($area in $scoreRanges.Areas) {
    $area.Value = 2
}

# Restore the view: scroll the frozen pane back to the top and select G3
$ws.Range("G3").Select() | Out-Null
